# Large code update: insert newly-found BvD IDs into the existing ID list
# in column A, keeping the existing entries (rows 1-6) untouched and
# re-writing rows 7 onward so the new IDs land next to the related
# existing ones, with the extra new rows appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "LULB166745",
    "DE2070071908",
    "DE8170085484",
    "US133277535L",
    "IT00079760328",
    "NL34275688",
    "FI01126502",
    "LULB176010",
    "LULB181081",
    "LULB185521",
    "LULB188712",
    "LULB185422",
    "HK0000244354",
    "US149146115L",
    "NL34140812"
)

$startRow = 7
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}
